$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3859567099616288

$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 0.4888122301146542

$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 0.5229064465527795

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 0.5446814307266112
